$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B ("Fixed") for rows that are missing it, and clear the stale
# error-message cells that are no longer needed now that the tests pass.

$ws.Range("B4").Value = "Fixed"
$ws.Range("D4").ClearContents()

$ws.Range("B11").Value = "Fixed"
$ws.Range("D11").ClearContents()

$ws.Range("B19").Value = "Fixed"

$ws.Range("B20").Value = "Fixed"
$ws.Range("E20").ClearContents()

$ws.Range("B21").Value = "Fixed"
$ws.Range("E21").ClearContents()

$ws.Range("B22").Value = "Fixed"
$ws.Range("E22").ClearContents()

$ws.Range("D11").Select()
